# "Generate Report for Archive"
#
# 1) The status text "Ready for handoff" becomes "In Translation" everywhere
#    it is used (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2).
# 2) The (now narrower) "Status"/zh-cn/de-de columns are resized down from
#    ~17.22 chars to ~13.41 chars: Overview columns E & F, and column C on
#    both the "zh-cn" and "de-de" sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1) Update the status text (same text, all occurrences, so they keep
#        sharing a single string in the saved workbook) ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2) Narrow the matching columns ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # column C
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # column C
